$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# A new run of 8.47 miles was logged on 2017-11-29 (row 8); extend the
# running total formula in I8 to include it, matching the author's
# "622 .... done !" commit.
$ws.Range("I8").Formula = "= 6.19 + 6.19 + 6.76 + 8.47"

$excel.Calculate()
